$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 325
$ws.Range("F7").Value = 874
$ws.Range("F8").Value = 58
$ws.Range("F9").Value = 522
$ws.Range("F10").Value = 67
$ws.Range("F11").Value = 297
$ws.Range("F12").Value = 1151
$ws.Range("F13").Value = 106
$ws.Range("F14").Value = 246
$ws.Range("F16").Value = 418
$ws.Range("F17").Value = 6667
$ws.Range("F21").Value = 7596
$ws.Range("F25").Value = 30
$ws.Range("F26").Value = 2113
$ws.Range("F27").Value = 898
$ws.Range("F28").Value = 4518
$ws.Range("F29").Value = 146
$ws.Range("F32").Value = 229
$ws.Range("F34").Value = 1716
$ws.Range("F36").Value = 178
$ws.Range("F39").Value = 1218
$ws.Range("F40").Value = 1805
$ws.Range("F41").Value = 2143
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 325
$ws.Range("F9").Value = 874
$ws.Range("F10").Value = 58
$ws.Range("F11").Value = 522
$ws.Range("F13").Value = 297
$ws.Range("F14").Value = 1151
$ws.Range("F16").Value = 106
$ws.Range("F17").Value = 246
$ws.Range("F19").Value = 418
$ws.Range("F20").Value = 6667
$ws.Range("F24").Value = 7596
$ws.Range("F28").Value = 30
$ws.Range("F29").Value = 2113
$ws.Range("F30").Value = 898
$ws.Range("F31").Value = 4518
$ws.Range("F32").Value = 146
$ws.Range("F36").Value = 229
$ws.Range("F38").Value = 1716
$ws.Range("F40").Value = 178
$ws.Range("F44").Value = 1218
$ws.Range("F45").Value = 1805
$ws.Range("F47").Value = 2143
